# Applies the changes described by the diff:
#   - WILMA P. BAYAS              -> PURISIMA CORAZON E. DUNGO   (x2)
#   - City Civil Registrar Office -> City Treasurer's Office     (x2)
#   - City Civil Registrar        -> Ticket Checker              (x2)
#   - March 01, 1984              -> August 01, 2013              (x2)
#   - August 31, 2009             -> December 31, 2022            (x2)
#   - "22nd" -> "23rd" in the issuance paragraph (superscript "rd" kept)
#
# Uses Range.Text (not Find.Execute's built-in replace) so that literal
# straight apostrophes are not auto-corrected into curly quotes.

$d = $word.ActiveDocument

function Replace-AllText($find, $replace) {
    $rng = $d.Content
    while ($rng.Find.Execute($find, $true, $false, $false, $false, $false, $true)) {
        $rng.Text = $replace
        $rng.Collapse(0) | Out-Null
        $moveCount = $d.Content.End - $rng.End
        if ($moveCount -gt 0) {
            $rng.MoveEnd(1, $moveCount) | Out-Null
        }
    }
}

# Order matters: replace the longer/more specific phrase before the shorter
# one that is a substring of it ("City Civil Registrar Office" contains
# "City Civil Registrar").
Replace-AllText "City Civil Registrar Office" "City Treasurer's Office"
Replace-AllText "City Civil Registrar" "Ticket Checker"
Replace-AllText "WILMA P. BAYAS" "PURISIMA CORAZON E. DUNGO"
Replace-AllText "March 01, 1984" "August 01, 2013"
Replace-AllText "August 31, 2009" "December 31, 2022"

# "This certification is issued this 22nd day of ..." -> "...23rd day of..."
# "22nd" is the only occurrence of that literal text in the document, and it
# spans a plain run ("22") followed by a superscript run ("nd"). Replace the
# two halves separately so the superscript formatting on "rd" is preserved.
$rng = $d.Content
if ($rng.Find.Execute("22nd", $true, $false, $false, $false, $false, $true)) {
    $start = $rng.Start
    $d.Range($start, $start + 2).Text = "23"
    $d.Range($start + 2, $start + 4).Text = "rd"
}
